$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Homme)
$ws.Range("B2").Value = 10528
$ws.Range("C2").Value = 83.1464223661349
$ws.Range("D2").Value = 28042
$ws.Range("E2").Value = 81.41330855882012
$ws.Range("F2").Value = 38570
$ws.Range("G2").Value = 81.879166135948708
$ws.Range("H2").Value = 5798
$ws.Range("I2").Value = 94.661224489795913

# Row 3 (Femme)
$ws.Range("B3").Value = 2134
$ws.Range("C3").Value = 16.853577633865111
$ws.Range("D3").Value = 6402
$ws.Range("E3").Value = 18.586691441179891
$ws.Range("F3").Value = 8536
$ws.Range("G3").Value = 18.120833864051288
$ws.Range("H3").Value = 327
$ws.Range("I3").Value = 5.3387755102040817

# Update selection to G17 as in the diff
$ws.Range("G17").Select()
